$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A92").Value = "2025-04-29 15:34:28"
$ws.Range("B92").Value = 228
